$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.841.21'
$ws.Range('E2').Value = '  -3.61%  '
$ws.Range('D3').Value = '3.310.91'
$ws.Range('E3').Value = '  -5.30%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '181.76'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -8.43%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '531.96'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.28%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.606'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.38%  '
$ws.Range('D8').Value = '3.308.52'
$ws.Range('E8').Value = '  -5.27%  '
$ws.Range('E9').Value = '  -0.06%  '
$ws.Range('E10').Value = '  -5.11%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '60.12'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.14%  '
$ws.Range('E12').Value = '  -5.93%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.16'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -6.46%  '
$ws.Range('D15').Value = '3.829.50'
$ws.Range('E15').Value = '  -5.73%  '
$ws.Range('D16').Value = '3.306.09'
$ws.Range('E16').Value = '  -5.49%  '
$ws.Range('E17').Value = '  -4.79%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '17.73'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.99%  '
$ws.Range('D19').Value = '64.673.13'
$ws.Range('E19').Value = '  -3.51%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.22'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.85%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.965'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.54%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '376.43'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.18%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.83'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.88%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '81.35'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.97%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.17'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -5.21%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.92'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.72%  '
$ws.Range('E27').Value = '  -0.97%  '
$ws.Range('E28').Value = '  -2.92%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '11.63'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.74%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.47'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.50%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '29.18'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -5.68%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '650.59'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.68%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.77'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.15%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.41'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.22%  '
$ws.Range('E35').Value = '  -3.34%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '59.20'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -6.25%  '
$ws.Range('E37').Value = '  +0.03%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.396'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.24%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '37.04'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.53%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.997'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.20%  '
$ws.Range('D41').Value = '0.0₃0704'
$ws.Range('E41').Value = '  +5.21%  '
$ws.Range('E42').Value = '  -1.58%  '
$ws.Range('D43').Value = '2.883.89'
$ws.Range('E43').Value = '  -6.31%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.51'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.31%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.72'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -8.91%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0403'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.01%  '
$ws.Range('E47').Value = '  -3.96%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.86'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +9.82%  '
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.128'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.10%  '
$ws.Range('B50').Value = 'ApeXProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.03'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.02%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.56'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.40%  '
